$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Actualizacion de la lista de bugs" - add bugs #17 and #18 to the BUGS list
# (rows 19 and 20), which previously were blank filler rows before the
# closing border row (now row 21).
# ---------------------------------------------------------------------------

# --- Row 19 : bug #17 -------------------------------------------------------
# Pick up the formatting that already exists elsewhere on the sheet for each
# column so no new style/font/fill entries need to be fabricated.
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B18").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null

$ws.Range("D18").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null

$ws.Range("F18").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Equipos para generar fixture"
$ws.Range("C19").Value = "Deberían ser como mínimo 2"
$ws.Range("D19").Value = "Flor"
$ws.Range("E19").Value = "edicion-configurar.aspx"
$ws.Range("F19").Value = "PENDIENTE"

# --- Row 20 : bug #18 -------------------------------------------------------
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null

$ws.Range("C17").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null

$ws.Range("D18").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null

$ws.Range("F18").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "modificacion de configuracion de edicion"
$ws.Range("C20").Value = "cuando se modific la congiuracion de la edicion deberia generar el fixture con los nuevos equipos"
$ws.Range("D20").Value = "Flor"
$ws.Range("E20").Value = "edicion-configurar.aspx"
$ws.Range("F20").Value = "PENDIENTE"

$ws.Rows.Item(20).RowHeight = 45

$excel.CutCopyMode = $false

# --- Cursor / viewport moved down to the newly edited rows -----------------
$ws.Application.Goto($ws.Range("C17"), $false)
$ws.Range("C17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
